# Fix duplicate-count bug and rebuild the "change from lang" breakdown.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Duplicate-count bug fix: 'en' count in Crossref should be 3436, not 3442 ---
$ws.Range("C8").Value2 = 3436

# --- 2. Move the grand-total formula for C2:C35 down from row 12 to row 15 ---
#        (value recalculates automatically once C8 is corrected: 5098)
$ws.Range("E12").ClearContents()
$ws.Range("E15").Formula = "=SUM(C2:C35)"

# --- 3. Split the single "change from lang..." bucket into three labelled sub-totals ---

# Drop the old single label + its grand total formula
$ws.Range("E37").ClearContents()
$ws.Range("E38").ClearContents()

# New label/total pair #1: overall "from one lang to another" (everything in the block minus the C41 subtotal row)
$ws.Range("E36").Value = "from one lang to another"
$ws.Range("G36").Formula = "=SUM(C36:C40,C42:C75)"

# New label/total pair #2: "en to another lang"
$ws.Range("E42").Value = "en to another lang"
$ws.Range("G42").Formula = "=SUM(C42:C52)"

# New label/total pair #3: "other lang to english"
$ws.Range("E45").Value = "other lang to english"
$ws.Range("G45").Formula = "=SUM(C36,C38,C39,C40,C53,C56,C57,C59,C60,C62,C63,C64,C65,C66,C67,C70,C71,C72,C73,C74,C75)"

# --- 4. View-state touch-ups to mirror the author's re-review of the sheet ---
$excel.ActiveWindow.ScrollRow = 38
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G46").Select()
$excel.ActiveWindow.Zoom = 100
